$d = $word.ActiveDocument

# --- Revert "add phone to export word" ---
# 1) Delete the extra "phone_approver" row (label + {phone_approver} value,
#    which also carried the stray _GoBack bookmark) that was appended
#    after the "approver" row in the 3rd table.
$t = $d.Tables(3)
$t.Rows($t.Rows.Count).Delete()

# 2) Move the _GoBack bookmark back onto the {phone} run (end of that
#    table cell's paragraph), where it lived before that commit.
$rng = $d.Content
$rng.Find.Execute("{phone}", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.InsertAfter("@@BOOKMARK_MARK@@")

$markRng = $d.Content
$markRng.Find.Execute("@@BOOKMARK_MARK@@", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$d.Bookmarks.Add("_GoBack", $markRng)
$markRng.Text = ""
